# Update "想去人数" (want-to-go count) figures on the 展览 (Exhibition) and
# 全部类型 (All Types) sheets to the freshly scraped values.
$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 3487
$ws1.Range("F5").Value = 2174
$ws1.Range("F6").Value = 424
$ws1.Range("F8").Value = 59
$ws1.Range("F10").Value = 1285
$ws1.Range("F11").Value = 235
$ws1.Range("F12").Value = 1696
$ws1.Range("F13").Value = 123

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 3487
$ws4.Range("F5").Value = 2174
$ws4.Range("F6").Value = 424
$ws4.Range("F9").Value = 59
$ws4.Range("F13").Value = 1285
$ws4.Range("F14").Value = 235
$ws4.Range("F15").Value = 1696
$ws4.Range("F16").Value = 123
